# Updating vocab for keyword extraction (ontology data)
#
# Rows that were previously flagged as "specific to Austrian market (where
# the exp was conducted)" / "specific to Estonian market (where the exp was
# conducted)" / "specific to exp conducted" (column H, comment) are
# reclassified:
#   - Rows whose tag itself already encodes the market (Austrian/Estonian
#     allergen names, "egg origin estonia") are now considered "No Change"
#     and their update_code (column G) flips from -1 to 0.
#   - A few remaining rows ("made in estonia", "country of origin of the
#     user") keep the market/experiment-specific comment, but the comment
#     text is re-cased to start with a capital letter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 27-39 (Austrian allergen tags) and 117-126 + 245 (Estonian tags /
# "country of the retailer") : update_code -> 0, comment -> "No Change"
$noChangeRows = @(27,28,29,30,31,32,33,34,35,36,37,38,39,117,118,119,120,121,122,123,124,125,126,245)
foreach ($r in $noChangeRows) {
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = "No Change"
}

# Rows 162-164 ("made in estonia") : comment re-cased to start with a
# capital "S" (update_code stays -1)
$estonianCapRows = @(162,163,164)
foreach ($r in $estonianCapRows) {
    $ws.Cells.Item($r, 8).Value = "Specific to Estonian market (where the exp was conducted)"
}

# Row 244 ("country of origin of the user") : comment re-cased to start
# with a capital "S" (update_code stays -1)
$ws.Cells.Item(244, 8).Value = "Specific to exp conducted"

# Restore the selection Excel had when the file was last saved
$ws.Range("H329").Select()
